$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 8 columns B through S from 20 to 2000
$ws.Range("B8:S8").Value = 2000

# Update selection to G12
$ws.Range("G12").Select()
